$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 16416586
$ws.Range("J43").Value = 25796908
$ws.Range("L43").Value = 25796908
$ws.Range("N43").Value = -25797046
$ws.Range("H70").Value = 1166.1666
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1166.1666
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 3498.4998
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -4038.4998
$ws.Range("H73").Value = 1166.1666
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1166.1666
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 3498.4998
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -5370.4998
$ws.Range("H80").Value = 803.05
$ws.Range("I80").Value = 1062.75
$ws.Range("J80").Value = 629.9167
$ws.Range("K80").Value = 3188.25
$ws.Range("L80").Value = 1889.7501
$ws.Range("M80").Value = -2190.25
$ws.Range("N80").Value = -3885.7501
$ws.Range("H83").Value = 803.05
$ws.Range("I83").Value = 1062.75
$ws.Range("J83").Value = 629.9167
$ws.Range("K83").Value = 9564.75
$ws.Range("L83").Value = 5669.2503
$ws.Range("M83").Value = -4572.75
$ws.Range("N83").Value = -15653.2503
$ws.Range("H86").Value = 3271.682
$ws.Range("I86").Value = 3544.8462
$ws.Range("J86").Value = 2877.111
$ws.Range("K86").Value = 3544.8462
$ws.Range("L86").Value = 2877.111
$ws.Range("M86").Value = -2421.8462
$ws.Range("N86").Value = -5123.111
$ws.Range("H89").Value = 3271.682
$ws.Range("I89").Value = 3544.8462
$ws.Range("J89").Value = 2877.111
$ws.Range("K89").Value = 17724.231
$ws.Range("L89").Value = 14385.555
$ws.Range("M89").Value = -12108.231
$ws.Range("N89").Value = -25617.555
$ws.Range("H94").Value = 3249
$ws.Range("I94").Value = 3249
$ws.Range("K94").Value = 3249
$ws.Range("M94").Value = -2798
$ws.Range("H112").Value = 2162
$ws.Range("J112").Value = 2294.75
$ws.Range("L112").Value = 6884.25
$ws.Range("N112").Value = -9100.25
$ws.Range("H116").Value = 2479.8438
$ws.Range("I116").Value = 2804
$ws.Range("J116").Value = 2310.0476
$ws.Range("K116").Value = 2804
$ws.Range("L116").Value = 2310.0476
$ws.Range("M116").Value = 638
$ws.Range("N116").Value = -9194.0476
$ws.Range("H129").Value = 880.35486
$ws.Range("J129").Value = 924.3043
$ws.Range("L129").Value = 2772.9129
$ws.Range("N129").Value = -12772.9129
$ws.Range("H132").Value = 9015647
$ws.Range("I132").Value = 9528950
$ws.Range("K132").Value = 28586850
$ws.Range("M132").Value = -28584320
$ws.Range("H135").Value = 1301.2778
$ws.Range("I135").Value = 580.3333
$ws.Range("J135").Value = 4906
$ws.Range("K135").Value = 5222.9997
$ws.Range("L135").Value = 44154
$ws.Range("M135").Value = -2687.9997
$ws.Range("N135").Value = -49224
$ws.Range("H138").Value = 2715.07
$ws.Range("J138").Value = 2695.1633
$ws.Range("L138").Value = 8085.4899
$ws.Range("N138").Value = -18365.4899

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 55556932
$ws.Range("I61").Value = 76923900
$ws.Range("J61").Value = 2799.8
$ws.Range("K61").Value = 76923900
$ws.Range("L61").Value = 2799.8
$ws.Range("M61").Value = -76923688
$ws.Range("N61").Value = -3223.8
$ws.Range("H74").Value = 2798
$ws.Range("I74").Value = 1597.1428
$ws.Range("J74").Value = 5600
$ws.Range("K74").Value = 1597.1428
$ws.Range("L74").Value = 5600
$ws.Range("M74").Value = -723.1428000000001
$ws.Range("N74").Value = -7348
$ws.Range("H77").Value = 2798
$ws.Range("I77").Value = 1597.1428
$ws.Range("J77").Value = 5600
$ws.Range("K77").Value = 7985.714
$ws.Range("L77").Value = 28000
$ws.Range("M77").Value = -3617.714
$ws.Range("N77").Value = -36736
$ws.Range("H132").Value = 3168.5278
$ws.Range("I132").Value = 2748.8845
$ws.Range("J132").Value = 4259.6
$ws.Range("K132").Value = 8246.6535
$ws.Range("L132").Value = 12778.8
$ws.Range("M132").Value = -5716.6535
$ws.Range("N132").Value = -17838.8
$ws.Range("H134").Value = 38040
$ws.Range("J134").Value = 38040
$ws.Range("L134").Value = 38040
$ws.Range("N134").Value = -48180
$ws.Range("H136").Value = 55556932
$ws.Range("I136").Value = 76923900
$ws.Range("J136").Value = 2799.8
$ws.Range("K136").Value = 230771700
$ws.Range("L136").Value = 8399.400000000001
$ws.Range("M136").Value = -230769150
$ws.Range("N136").Value = -13499.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 36087.6
$ws.Range("J58").Value = 36087.6
$ws.Range("L58").Value = 36087.6
$ws.Range("N58").Value = -36675.6
$ws.Range("H134").Value = 5812.8184
$ws.Range("I134").Value = 850.7857
$ws.Range("K134").Value = 2552.3571
$ws.Range("M134").Value = -17.35710000000017

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1637.6825
$ws.Range("I31").Value = 1598.305
$ws.Range("J31").Value = 2218.5
$ws.Range("K31").Value = 1598.305
$ws.Range("L31").Value = 2218.5
$ws.Range("M31").Value = -1303.305
$ws.Range("N31").Value = -2808.5
$ws.Range("H34").Value = 1637.6825
$ws.Range("I34").Value = 1598.305
$ws.Range("J34").Value = 2218.5
$ws.Range("K34").Value = 1598.305
$ws.Range("L34").Value = 2218.5
$ws.Range("M34").Value = -1396.305
$ws.Range("N34").Value = -2622.5
$ws.Range("H43").Value = 15149.111
$ws.Range("J43").Value = 15149.111
$ws.Range("L43").Value = 15149.111
$ws.Range("N43").Value = -15517.111
$ws.Range("H58").Value = 5576.9165
$ws.Range("I58").Value = 857.61536
$ws.Range("J58").Value = 11154.272
$ws.Range("K58").Value = 857.61536
$ws.Range("L58").Value = 11154.272
$ws.Range("M58").Value = -654.61536
$ws.Range("N58").Value = -11560.272
$ws.Range("H101").Value = 15149.111
$ws.Range("J101").Value = 15149.111
$ws.Range("L101").Value = 15149.111
$ws.Range("N101").Value = -21639.111
$ws.Range("H136").Value = 5576.9165
$ws.Range("I136").Value = 857.61536
$ws.Range("J136").Value = 11154.272
$ws.Range("K136").Value = 2572.84608
$ws.Range("L136").Value = 33462.81600000001
$ws.Range("M136").Value = -22.84608000000026
$ws.Range("N136").Value = -38562.81600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2327.2727
$ws.Range("J55").Value = 2937.5
$ws.Range("L55").Value = 8812.5
$ws.Range("N55").Value = -9166.5
$ws.Range("H114").Value = 496.55554
$ws.Range("I114").Value = 295.7
$ws.Range("J114").Value = 747.625
$ws.Range("K114").Value = 887.0999999999999
$ws.Range("L114").Value = 2242.875
$ws.Range("M114").Value = 2366.9
$ws.Range("N114").Value = -8750.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 27793.625
$ws.Range("J106").Value = 27793.625
$ws.Range("L106").Value = 27793.625
$ws.Range("N106").Value = -30317.625
$ws.Range("H136").Value = 2240.0667
$ws.Range("I136").Value = 2185.7856
$ws.Range("K136").Value = 6557.3568
$ws.Range("M136").Value = -4007.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H97").Value = 17535
$ws.Range("J97").Value = 17535
$ws.Range("L97").Value = 17535
$ws.Range("N97").Value = -19517
$ws.Range("H101").Value = 15080
$ws.Range("J101").Value = 15080
$ws.Range("L101").Value = 15080
$ws.Range("N101").Value = -21570
$ws.Range("H104").Value = 9799.5
$ws.Range("J104").Value = 9799.5
$ws.Range("L104").Value = 9799.5
$ws.Range("N104").Value = -16787.5
$ws.Range("H126").Value = 125001340
$ws.Range("I126").Value = 250001820
$ws.Range("J126").Value = 860
$ws.Range("K126").Value = 750005460
$ws.Range("L126").Value = 2580
$ws.Range("M126").Value = -750002990
$ws.Range("N126").Value = -7520
$ws.Range("H136").Value = 1187.295
$ws.Range("I136").Value = 446.92307
$ws.Range("J136").Value = 2499.7727
$ws.Range("K136").Value = 1340.76921
$ws.Range("L136").Value = 7499.3181
$ws.Range("M136").Value = 1209.23079
$ws.Range("N136").Value = -12599.3181
$ws.Range("H137").Value = 37728
$ws.Range("J137").Value = 37728
$ws.Range("L137").Value = 37728
$ws.Range("N137").Value = -47928
